$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in column F ("blogspot") that need to be filled in with "ok"
$rows = @(4,5,6,7,9,10,11,12,13,14,15,17)
foreach ($r in $rows) {
    $ws.Range("F$r").Value = "ok"
}

# Move the active selection to A19 (below the table), matching the saved view state
$ws.Range("A19").Select()
